$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a data row (row 6) for CARLA, plus a blank spacer row (row 7),
# pushing the "Documents" block (previously rows 7-10) down to rows 9-12.
$ws.Rows("6:7").Insert()

# Insert one more blank row before the "Understanding Geotab" block
# (previously row 12, now at row 14 after the first insert) so it lands
# on row 15, matching the target layout.
$ws.Rows("14:14").Insert()

# B6 gets a hyperlink to the CARLA playlist, with the auto-generated
# lowercase ScreenTip that Excel produces for a typed-in URL.
$url = "https://www.youtube.com/playlist?list=PLQVvvaa0QuDeI12McNQdnTlWz9XlCa0uo"
$tip = "https://www.youtube.com/playlist?list=plqvvvaa0qudei12mcnqdntlwz9xlca0uo"
$ws.Hyperlinks.Add($ws.Range("B6"), $url, [Type]::Missing, $tip)
$ws.Range("B6").Style = $ws.Range("B2").Style

# Fill in the rest of the new CARLA row.
$ws.Range("A6").Value = "CARLA (First 2 videos)"
$ws.Range("C6").Value = "1:00 Hrs"
$ws.Range("D6").Value = "Optional"

# Give the blank spacer row B7 the same Hyperlink-ish formatting left
# behind from the copied row (no value, just the style).
$ws.Range("B7").Style = $ws.Range("B2").Style

# Match the saved selection state from the authored workbook.
$ws.Range("B4").Select()
